$d = $word.ActiveDocument

$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$endRange.Collapse(0)
$endRange.InsertAfter("Меня зовут Настя")
